$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add Job (C) and City (D) columns, and update Row 5's age value.
$ws.Range("C1").Value = "エンジニア"
$ws.Range("D1").Value = "東京"

$ws.Range("C2").Value = "デザイナー"
$ws.Range("D2").Value = "大阪"

$ws.Range("C3").Value = "営業"
$ws.Range("D3").Value = "福岡"

$ws.Range("C4").Value = "マネージャー"
$ws.Range("D4").Value = "名古屋"

$ws.Range("B5").Value = 22.0
$ws.Range("C5").Value = "デザイナー"
$ws.Range("D5").Value = "横浜"
